# Apply "output generated at 456a3b4" updates to the 上海-漫展信息 workbook.
# Each sheet gets a handful of "想去人数" (F column) bumps, and one listing
# on 演出 (Performance) flips its "最低票价" (G column) from a numeric price
# to the "已停售" (sales stopped) status string.

$wb = $excel.ActiveWorkbook

# ---- Sheet 展览 (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 640
$ws1.Range("F6").Value = 2824
$ws1.Range("F12").Value = 327
$ws1.Range("F14").Value = 5981
$ws1.Range("F17").Value = 12
$ws1.Range("F18").Value = 236
$ws1.Range("F21").Value = 543
$ws1.Range("F22").Value = 33
$ws1.Range("F23").Value = 31
$ws1.Range("F24").Value = 107
$ws1.Range("F25").Value = 1319
$ws1.Range("F28").Value = 48
$ws1.Range("F29").Value = 2061
$ws1.Range("F30").Value = 183
$ws1.Range("F31").Value = 357
$ws1.Range("F33").Value = 3305

# ---- Sheet 演出 (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G5").Value = "已停售"
$ws2.Range("F12").Value = 646
$ws2.Range("F15").Value = 1006
$ws2.Range("F17").Value = 76
$ws2.Range("F22").Value = 352
$ws2.Range("F24").Value = 4055
$ws2.Range("F26").Value = 13
$ws2.Range("F28").Value = 144

# ---- Sheet 本地生活 (Local Life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2600
$ws3.Range("F6").Value = 1143
$ws3.Range("F12").Value = 646

# ---- Sheet 全部类型 (All Types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2600
$ws4.Range("F6").Value = 1143
$ws4.Range("F11").Value = 640
$ws4.Range("F12").Value = 2824
$ws4.Range("F14").Value = 646
$ws4.Range("F18").Value = 327
$ws4.Range("F20").Value = 5981
$ws4.Range("F24").Value = 236
$ws4.Range("F27").Value = 543
$ws4.Range("F29").Value = 76
$ws4.Range("F32").Value = 31
$ws4.Range("F34").Value = 352
$ws4.Range("F37").Value = 13
$ws4.Range("F39").Value = 144
$ws4.Range("F41").Value = 48
$ws4.Range("F44").Value = 2061
$ws4.Range("F47").Value = 183
$ws4.Range("F48").Value = 357
$ws4.Range("F50").Value = 3305
